$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B5 value from 602 to 604
$ws.Range("B5").Value = 604

# Add new row 6: semana 5, casos 407
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 407

# Add new row 7: semana 6, casos 13
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 13
